# Insert two new price-report rows (Murcott, Provincia de Limarí, reported
# 2021-09-13 / serial 44460) right after the existing row 90, pushing all
# subsequent rows (old 91..127) down by two. This matches the commit's
# weekly Fruta/Hortaliza refresh: new rows 91-92 are added, and the sheet
# grows from A1:T127 to A1:T129.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91 downward (and everything below it) by inserting 2 new blank
# rows at position 91.
$ws.Range("91:92").EntireRow.Insert()

# New row 91: Murcott / Primera
$ws.Cells.Item(91, 1).Value  = 7
$ws.Cells.Item(91, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(91, 3).Value  = "Ñuble"
$ws.Cells.Item(91, 4).Value  = 44460
$ws.Cells.Item(91, 5).Value  = 16
$ws.Cells.Item(91, 6).Value  = "Fruta"
$ws.Cells.Item(91, 7).Value  = 100102
$ws.Cells.Item(91, 8).Value  = "Cítricos"
$ws.Cells.Item(91, 9).Value  = 100102004
$ws.Cells.Item(91, 10).Value = "Mandarina"
$ws.Cells.Item(91, 11).Value = "Murcott"
$ws.Cells.Item(91, 12).Value = "Primera"
$ws.Cells.Item(91, 13).Value = 300
$ws.Cells.Item(91, 14).Value = 5500
$ws.Cells.Item(91, 15).Value = 6000
$ws.Cells.Item(91, 16).Value = 5750
$ws.Cells.Item(91, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(91, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(91, 19).Value = 575
$ws.Cells.Item(91, 20).Value = 10

# New row 92: Murcott / Segunda
$ws.Cells.Item(92, 1).Value  = 7
$ws.Cells.Item(92, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(92, 3).Value  = "Ñuble"
$ws.Cells.Item(92, 4).Value  = 44460
$ws.Cells.Item(92, 5).Value  = 16
$ws.Cells.Item(92, 6).Value  = "Fruta"
$ws.Cells.Item(92, 7).Value  = 100102
$ws.Cells.Item(92, 8).Value  = "Cítricos"
$ws.Cells.Item(92, 9).Value  = 100102004
$ws.Cells.Item(92, 10).Value = "Mandarina"
$ws.Cells.Item(92, 11).Value = "Murcott"
$ws.Cells.Item(92, 12).Value = "Segunda"
$ws.Cells.Item(92, 13).Value = 90
$ws.Cells.Item(92, 14).Value = 5000
$ws.Cells.Item(92, 15).Value = 5000
$ws.Cells.Item(92, 16).Value = 5000
$ws.Cells.Item(92, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(92, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(92, 19).Value = 500
$ws.Cells.Item(92, 20).Value = 10
